# update dac diem khoa
#
# The source edits in the commit only merge adjacent <w:r> runs that sit
# next to each other (same formatting) into a single run -- i.e. the
# visible text does not change, only how it is split across runs -- plus
# one paragraph gets a brand-new trailing run consisting of four spaces.
#
# We implement every change with the same primitive: locate the exact
# text with Find, then replace that found range with a hand-built OOXML
# fragment (via Range.InsertXML) describing precisely the run layout we
# want. This lets us control run boundaries exactly (merging some runs,
# keeping a leading <w:tab/> untouched, or adding a brand new run)
# without relying on however Find/Replace happens to split runs.

function New-WordBodyXml($innerP) {
    return @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
$innerP
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
}

function Replace-FoundRange($d, $searchText, $innerP) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    # Re-seat the hit in a fresh, independent Range before mutating --
    # InsertXML-ing the live Find range directly corrupts surrounding text.
    $target = $d.Range($rng.Start, $rng.End)
    $xml = New-WordBodyXml $innerP
    $target.InsertXML($xml)
}

$d = $word.ActiveDocument

# 1) "Mở khóa bằng nhận diện gương mặt: " + "tăng khả năng bảo mật, có thể
#    kết hợp với lại " -> single run.
Replace-FoundRange $d "Mở khóa bằng nhận diện gương mặt: tăng khả năng bảo mật, có thể kết hợp với lại " `
    '<w:p><w:r><w:t xml:space="preserve">Mở khóa bằng nhận diện gương mặt: tăng khả năng bảo mật, có thể kết hợp với lại </w:t></w:r></w:p>'

# 2) tab run's lone space " " + "+ Đối tượng sử dụng: ..." -> merge the
#    two <w:t> into the run that already carries the <w:tab/>. Including
#    the literal tab char at the front of the search text makes Find's
#    hit start exactly at the <w:tab/>, so the whole original run (tab +
#    space) is cleanly replaced instead of leaving a stray remnant run.
Replace-FoundRange $d "`t + Đối tượng sử dụng: người lớn hoặc trẻ em (> 1m3) (vì module camera được đặt cao khoảng 1m3)" `
    '<w:p><w:r><w:tab/><w:t xml:space="preserve"> + Đối tượng sử dụng: người lớn hoặc trẻ em (&gt; 1m3) (vì module camera được đặt cao khoảng 1m3)</w:t></w:r></w:p>'

# 3) tab run's "+ " + "Lưu được bao nhiêu gương mặt? Lưu ở đâu?" -> merge,
#    keeping the leading <w:tab/> (same literal-tab-in-search trick as #2).
Replace-FoundRange $d "`t+ Lưu được bao nhiêu gương mặt? Lưu ở đâu?" `
    '<w:p><w:r><w:tab/><w:t>+ Lưu được bao nhiêu gương mặt? Lưu ở đâu?</w:t></w:r></w:p>'

# 4) "-" + " Website: gồm những chức năng:" -> single run.
Replace-FoundRange $d "- Website: gồm những chức năng:" `
    '<w:p><w:r><w:t>- Website: gồm những chức năng:</w:t></w:r></w:p>'

# 5) tab run's "+ " + "Mở cửa từ xa" -> merge (trailing lone-space run
#    after it is untouched; same literal-tab-in-search trick as #2).
Replace-FoundRange $d "`t+ Mở cửa từ xa" `
    '<w:p><w:r><w:tab/><w:t>+ Mở cửa từ xa</w:t></w:r></w:p>'

# 6) "- " + "Cảnh báo sms: " + "khi nhận dạng sai 5 lần..." -> single run.
Replace-FoundRange $d "- Cảnh báo sms: khi nhận dạng sai 5 lần đối với gương mặt hoặc vân tay, sai 3 lần đối với mặt khẩu số." `
    '<w:p><w:r><w:t>- Cảnh báo sms: khi nhận dạng sai 5 lần đối với gương mặt hoặc vân tay, sai 3 lần đối với mặt khẩu số.</w:t></w:r></w:p>'

# 7) "- Tích hợp màn hình ... là ai (" + "khi có khách ... trong nhà" + ")"
#    -> single run.
Replace-FoundRange $d "- Tích hợp màn hình hiển thị bên trong cửa để xem người bên ngoài nhấn chuông là ai (khi có khách nhấn chuông, thì camera chuyển sang chế độ streaming và streaming mặt khách vô trong nhà)" `
    '<w:p><w:r><w:t>- Tích hợp màn hình hiển thị bên trong cửa để xem người bên ngoài nhấn chuông là ai (khi có khách nhấn chuông, thì camera chuyển sang chế độ streaming và streaming mặt khách vô trong nhà)</w:t></w:r></w:p>'

# 8) Append a brand new trailing run of four spaces after the existing
#    run in the "Tích hợp chìa khóa cơ ..." paragraph.
Replace-FoundRange $d "- Tích hợp chìa khóa cơ để mở cửa trường hợp nhà mất điện." `
    '<w:p><w:r><w:t>- Tích hợp chìa khóa cơ để mở cửa trường hợp nhà mất điện.</w:t></w:r><w:r><w:t xml:space="preserve">    </w:t></w:r></w:p>'
